$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 5.922420666666667
$ws.Range("H2").Value = 17.767262
$ws.Range("I2").Value = 0.5833698282960434
$ws.Range("J2").Value = 0.6311054116979437
$ws.Range("M2").Value = 7.004922666666666
$ws.Range("N2").Value = 21.014768
$ws.Range("O2").Value = 0.134029393318039
$ws.Range("P2").Value = 0.1464771065395205
$ws.Range("Q2").Value = 41.48609876946844
$ws.Range("R2").Value = 373.374888925216
$ws.Range("S2").Value = 0.07818870416656726
$ws.Range("T2").Value = 0.09244249462694763
$ws.Range("G3").Value = 5.922420666666667
$ws.Range("H3").Value = 17.767262
$ws.Range("I3").Value = 0.5833698282960434
$ws.Range("J3").Value = 0.6311054116979437
$ws.Range("O3").Value = 0.3796257919253833
$ws.Range("P3").Value = 0.4148827819958515
$ws.Range("Q3").Value = 117.5055165838347
$ws.Range("R3").Value = 1057.549649254512
$ws.Range("S3").Value = 0.2214622330522603
$ws.Range("T3").Value = 0.2618347689378801
$ws.Range("G4").Value = 5.922420666666667
$ws.Range("H4").Value = 17.767262
$ws.Range("I4").Value = 0.5833698282960434
$ws.Range("J4").Value = 0.6311054116979437
$ws.Range("M4").Value = 6.293636666666667
$ws.Range("N4").Value = 18.88091
$ws.Range("O4").Value = 0.12041993100245
$ws.Range("P4").Value = 0.1316036924905903
$ws.Range("Q4").Value = 37.27356386315778
$ws.Range("R4").Value = 335.46207476842
$ws.Range("S4").Value = 0.07024935447232068
$ws.Range("T4").Value = 0.08305580253024357
$ws.Range("G5").Value = 5.922420666666667
$ws.Range("H5").Value = 17.767262
$ws.Range("I5").Value = 0.5833698282960434
$ws.Range("J5").Value = 0.6311054116979437
$ws.Range("M5").Value = 13.324299
$ws.Range("N5").Value = 26.648598
$ws.Range("O5").Value = 0.2549418168249328
$ws.Range("P5").Value = 0.1857460205306503
$ws.Range("Q5").Value = 78.91210376644599
$ws.Range("R5").Value = 473.472622598676
$ws.Range("S5").Value = 0.1487253639066424
$ws.Range("T5").Value = 0.1172253187582507
$ws.Range("G6").Value = 5.922420666666667
$ws.Range("H6").Value = 17.767262
$ws.Range("I6").Value = 0.5833698282960434
$ws.Range("J6").Value = 0.6311054116979437
$ws.Range("M6").Value = 5.800427666666667
$ws.Range("N6").Value = 17.401283
$ws.Range("O6").Value = 0.110983066929195
$ws.Range("P6").Value = 0.1212903984433873
$ws.Range("Q6").Value = 34.35257268857178
$ws.Range("R6").Value = 309.1731541971459
$ws.Range("S6").Value = 0.06474417269825278
$ws.Range("T6").Value = 0.0765470268446216
$ws.Range("I7").Value = 0.1897160182974547
$ws.Range("J7").Value = 0.2052399696141807
$ws.Range("M7").Value = 7.004922666666666
$ws.Range("N7").Value = 21.014768
$ws.Range("O7").Value = 0.134029393318039
$ws.Range("P7").Value = 0.1464771065395205
$ws.Range("Q7").Value = 13.49157445496889
$ws.Range("R7").Value = 121.42417009472
$ws.Range("S7").Value = 0.02542752283512183
$ws.Range("T7").Value = 0.03006295689534429
$ws.Range("I8").Value = 0.1897160182974547
$ws.Range("J8").Value = 0.2052399696141807
$ws.Range("O8").Value = 0.3796257919253833
$ws.Range("P8").Value = 0.4148827819958515
$ws.Range("S8").Value = 0.07202109368710173
$ws.Range("T8").Value = 0.08515052957027533
$ws.Range("I9").Value = 0.1897160182974547
$ws.Range("J9").Value = 0.2052399696141807
$ws.Range("M9").Value = 6.293636666666667
$ws.Range("N9").Value = 18.88091
$ws.Range("O9").Value = 0.12041993100245
$ws.Range("P9").Value = 0.1316036924905903
$ws.Range("Q9").Value = 12.12162813515556
$ws.Range("R9").Value = 109.0946532164
$ws.Range("S9").Value = 0.02284558983343904
$ws.Range("T9").Value = 0.02701033784788274
$ws.Range("I10").Value = 0.1897160182974547
$ws.Range("J10").Value = 0.2052399696141807
$ws.Range("M10").Value = 13.324299
$ws.Range("N10").Value = 26.648598
$ws.Range("O10").Value = 0.2549418168249328
$ws.Range("P10").Value = 0.1857460205306503
$ws.Range("Q10").Value = 25.66277753132
$ws.Range("R10").Value = 153.97666518792
$ws.Range("S10").Value = 0.0483665463855453
$ws.Range("T10").Value = 0.03812250760966564
$ws.Range("I11").Value = 0.1897160182974547
$ws.Range("J11").Value = 0.2052399696141807
$ws.Range("M11").Value = 5.800427666666667
$ws.Range("N11").Value = 17.401283
$ws.Range("O11").Value = 0.110983066929195
$ws.Range("P11").Value = 0.1212903984433873
$ws.Range("Q11").Value = 11.17170102503556
$ws.Range("R11").Value = 100.54530922532
$ws.Range("S11").Value = 0.02105526555624679
$ws.Range("T11").Value = 0.02489363769101269
$ws.Range("G12").Value = 2.303652
$ws.Range("H12").Value = 4.607303999999999
$ws.Range("I12").Value = 0.2269141534065018
$ws.Range("J12").Value = 0.1636546186878756
$ws.Range("M12").Value = 7.004922666666666
$ws.Range("N12").Value = 21.014768
$ws.Range("O12").Value = 0.134029393318039
$ws.Range("P12").Value = 0.1464771065395205
$ws.Range("Q12").Value = 16.136904110912
$ws.Range("R12").Value = 96.82142466547198
$ws.Range("S12").Value = 0.03041316631634986
$ws.Range("T12").Value = 0.02397165501722856
$ws.Range("G13").Value = 2.303652
$ws.Range("H13").Value = 4.607303999999999
$ws.Range("I13").Value = 0.2269141534065018
$ws.Range("J13").Value = 0.1636546186878756
$ws.Range("O13").Value = 0.3796257919253833
$ws.Range("P13").Value = 0.4148827819958515
$ws.Range("Q13").Value = 45.70628017238399
$ws.Range("R13").Value = 274.2376810343039
$ws.Range("S13").Value = 0.08614246518602116
$ws.Range("T13").Value = 0.06789748348769611
$ws.Range("G14").Value = 2.303652
$ws.Range("H14").Value = 4.607303999999999
$ws.Range("I14").Value = 0.2269141534065018
$ws.Range("J14").Value = 0.1636546186878756
$ws.Range("M14").Value = 6.293636666666667
$ws.Range("N14").Value = 18.88091
$ws.Range("O14").Value = 0.12041993100245
$ws.Range("P14").Value = 0.1316036924905903
$ws.Range("Q14").Value = 14.49834869444
$ws.Range("R14").Value = 86.99009216663998
$ws.Range("S14").Value = 0.02732498669669031
$ws.Range("T14").Value = 0.021537552112464
$ws.Range("G15").Value = 2.303652
$ws.Range("H15").Value = 4.607303999999999
$ws.Range("I15").Value = 0.2269141534065018
$ws.Range("J15").Value = 0.1636546186878756
$ws.Range("M15").Value = 13.324299
$ws.Range("N15").Value = 26.648598
$ws.Range("O15").Value = 0.2549418168249328
$ws.Range("P15").Value = 0.1857460205306503
$ws.Range("Q15").Value = 30.69454803994799
$ws.Range("R15").Value = 122.778192159792
$ws.Range("S15").Value = 0.05784990653274509
$ws.Range("T15").Value = 0.03039819416273389
$ws.Range("G16").Value = 2.303652
$ws.Range("H16").Value = 4.607303999999999
$ws.Range("I16").Value = 0.2269141534065018
$ws.Range("J16").Value = 0.1636546186878756
$ws.Range("M16").Value = 5.800427666666667
$ws.Range("N16").Value = 17.401283
$ws.Range("O16").Value = 0.110983066929195
$ws.Range("P16").Value = 0.1212903984433873
$ws.Range("Q16").Value = 13.362166795172
$ws.Range("R16").Value = 80.17300077103198
$ws.Range("S16").Value = 0.02518362867469541
$ws.Range("T16").Value = 0.01984973390775306
